$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column I ("Target Value") numeric updates -------------------------------
$ws.Range("I3").Value  = 5
$ws.Range("I4").Value  = 69
$ws.Range("I5").Value  = 22
$ws.Range("I6").Value  = 2
$ws.Range("I9").Value  = 2
$ws.Range("I10").Value = 96
$ws.Range("I13").Value = 97
$ws.Range("I16").Value = 24
$ws.Range("I17").Value = 29
$ws.Range("I18").Value = 38
$ws.Range("I20").Value = 7
$ws.Range("I22").Value = 40
$ws.Range("I23").Value = 44
$ws.Range("I24").Value = 13
$ws.Range("I28").Value = 40
$ws.Range("I29").Value = 44
$ws.Range("I30").Value = 13
$ws.Range("I34").Value = 40
$ws.Range("I35").Value = 44
$ws.Range("I36").Value = 13
$ws.Range("I45").Value = 97
$ws.Range("I52").Value = 78
$ws.Range("I53").Value = 4
$ws.Range("I54").Value = 16

# --- Row 14: previously a blank "Check ResidualCapacity" placeholder row,    --
# --- now filled in like the other tech rows (target 2050, value 1, etc.)    --
$ws.Range("G14").Value = "Check ResidualCapacity"
$ws.Range("H14").Value = 2050
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = "Percent"
$ws.Range("K14").Value = "Lower"
$ws.Range("L14").Value = "Linear"
$ws.Range("M14").Value = "Continuous"

# --- Row 25: same pattern as row 14, but the target year is 2051 -------------
$ws.Range("G25").Value = "Check ResidualCapacity"
$ws.Range("H25").Value = 2051
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = "Percent"
$ws.Range("K25").Value = "Lower"
$ws.Range("L25").Value = "Linear"
$ws.Range("M25").Value = "Continuous"

# --- Row 47: the previously filled-in row is cleared back to blank ----------
$ws.Range("G47:M47").ClearContents()
